# Corecion a Diebold Mariano: insert a "d=6" row into the results table.
#
# The table lists rows d=1, d=2, d=3, d=4, d=5, d=7, d=10 (column A) with
# four result columns (B:E). A new row for d=6 needs to be inserted right
# before the existing d=7 row, pushing d=7 and d=10 down by one row, and
# d=7's row gets freshly recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old layout: row7 = d=7, row8 = d=10.
# Insert a blank row at position 7 -> row7 becomes blank, old row7 (d=7)
# shifts to row8, old row8 (d=10) shifts to row9.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the d=6 label and its values.
$ws.Range("A7").Value = "d=6"
$ws.Range("B7").Value = 97.98222855450359
$ws.Range("C7").Value = 98.01050164382112
$ws.Range("D7").Value = 98.10180469708054
$ws.Range("E7").Value = 98.05001601071173

# Insert() doesn't carry over the thin-border box style used by the other
# row-label cells in column A, so restore it on the new label cell.
$ws.Range("A7").Borders.Item(1).LineStyle = 1
$ws.Range("A7").Borders.Item(2).LineStyle = 1
$ws.Range("A7").Borders.Item(3).LineStyle = 1
$ws.Range("A7").Borders.Item(4).LineStyle = 1
